# Actualización desde MV -datos-
# Refresh the latest quarter (2021-Q1) with revised figures and append the
# newly published quarter (2021-Q2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 (period 01-01-2021) with refreshed figures ---
$ws.Range("B74").Value = 50544
$ws.Range("C74").Value = -3004
$ws.Range("D74").Value = 5754
$ws.Range("E74").Value = -11
$ws.Range("F74").Value = 53283
$ws.Range("G74").Value = 41167
$ws.Range("H74").Value = 12116
$ws.Range("I74").Value = 13449
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = -1332

# --- Append new row 75 (period 01-04-2021) ---
# Write the period label via a scratch formula cell + copy/paste-special so
# the text is stored as a shared string (matching the existing "Serie"
# column) instead of being auto-converted to a date serial number.
$ws.Range("Z1").Formula = "=""01-04-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("B75").Value = 51330
$ws.Range("C75").Value = -3387
$ws.Range("D75").Value = 5963
$ws.Range("E75").Value = 35
$ws.Range("F75").Value = 53940
$ws.Range("G75").Value = 43315
$ws.Range("H75").Value = 10625
$ws.Range("I75").Value = 12521
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = -1895
